# Add 8 new flight-arrival rows (Friday, Jan 13) to the "Main Data" sheet,
# directly below the existing data (which ends at row 51).
#
# Columns: A=NUMBER, B=DATE, C=TIME, D=FLIGHT, E=FROM, F=SHORT, G=AIRLINE,
#          H=MODEL, I=AIRCFAT ID, J=STATUS, K=(unused), L=DIFFERENCE, M=(unused)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @(51, "Friday, Jan 13", "5:49 AM",  "X7592",  "New York", "(JFK)", "Challenge Airlines ", "B744", "(OO-ACE)",  "5:39 AM",  "0 hours, -10 minutes"),
    @(52, "Friday, Jan 13", "7:30 AM",  "LH9931", "Munich",   "(MUC)", "Lufthansa ",           "A319", "(D-AILW)",  "7:27 AM",  "0 hours, -3 minutes"),
    @(53, "Friday, Jan 13", "10:15 AM", "SK7181", "Oslo",     "(OSL)", "SAS ",                 "B737", "(LN-RPJ)",  "10:17 AM", "0 hours, 2 minutes"),
    @(54, "Friday, Jan 13", "11:17 AM", "5Y309",  "Dover",    "(DOV)", "Atlas Air ",           "B744", "(N452PA)",  "11:18 AM", "0 hours, 1 minutes"),
    @(55, "Friday, Jan 13", "11:25 AM", "LO3809", "Warsaw",   "(WAW)", "LOT ",                 "E195", "(SP-LNK)",  "11:10 AM", "0 hours, -15 minutes"),
    @(56, "Friday, Jan 13", "1:10 PM",  "FR8224", "Bristol",  "(BRS)", "Ryanair ",             "B738", "(EI-DCM)",  "1:03 PM",  "0 hours, -7 minutes"),
    @(57, "Friday, Jan 13", "2:07 PM",  "UNKNOWN","Belgrade", "(BEG)", "AirPink ",             "C525", "(YU-MTU)",  "2:19 PM",  "0 hours, 12 minutes"),
    @(58, "Friday, Jan 13", "2:15 PM",  "LO3801", "Warsaw",   "(WAW)", "LOT ",                 "E190", "(SP-LMH)",  "2:06 PM",  "0 hours, -9 minutes")
)

$startRow = 52
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $data = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $data[0]    # A - NUMBER
    $ws.Cells.Item($r, 2).Value = $data[1]    # B - DATE
    $ws.Cells.Item($r, 3).Value = $data[2]    # C - TIME
    $ws.Cells.Item($r, 4).Value = $data[3]    # D - FLIGHT
    $ws.Cells.Item($r, 5).Value = $data[4]    # E - FROM
    $ws.Cells.Item($r, 6).Value = $data[5]    # F - SHORT
    $ws.Cells.Item($r, 7).Value = $data[6]    # G - AIRLINE
    $ws.Cells.Item($r, 8).Value = $data[7]    # H - MODEL
    $ws.Cells.Item($r, 9).Value = $data[8]    # I - AIRCFAT ID
    $ws.Cells.Item($r, 10).Value = $data[9]   # J - STATUS
    $ws.Cells.Item($r, 12).Value = $data[10]  # L - DIFFERENCE
}
